# Update "想去人数" (interest count) figures in column F across all four
# sheets of the 广州-漫展信息 workbook, reflecting freshly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 14098
$ws.Range("F3").Value = 13805
$ws.Range("F4").Value = 832
$ws.Range("F9").Value = 74
$ws.Range("F10").Value = 797
$ws.Range("F11").Value = 2182
$ws.Range("F12").Value = 179
$ws.Range("F13").Value = 124
$ws.Range("F14").Value = 98
$ws.Range("F15").Value = 224
$ws.Range("F17").Value = 594
$ws.Range("F19").Value = 511
$ws.Range("F20").Value = 345
$ws.Range("F21").Value = 32
$ws.Range("F22").Value = 305
$ws.Range("F23").Value = 889
$ws.Range("F24").Value = 148
$ws.Range("F25").Value = 73
$ws.Range("F26").Value = 28
$ws.Range("F29").Value = 79
$ws.Range("F30").Value = 37

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 55
$ws.Range("F6").Value = 115
$ws.Range("F7").Value = 185
$ws.Range("F8").Value = 1871
$ws.Range("F13").Value = 86
$ws.Range("F15").Value = 1842

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 234
$ws.Range("F3").Value = 163

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 234
$ws.Range("F3").Value = 14098
$ws.Range("F4").Value = 13805
$ws.Range("F5").Value = 832
$ws.Range("F10").Value = 74
$ws.Range("F11").Value = 797
$ws.Range("F12").Value = 55
$ws.Range("F14").Value = 2182
$ws.Range("F15").Value = 163
$ws.Range("F16").Value = 179
$ws.Range("F17").Value = 179
$ws.Range("F18").Value = 124
$ws.Range("F19").Value = 98
$ws.Range("F20").Value = 224
$ws.Range("F24").Value = 115
$ws.Range("F26").Value = 594
$ws.Range("F28").Value = 511
$ws.Range("F29").Value = 345
$ws.Range("F30").Value = 32
$ws.Range("F31").Value = 305
$ws.Range("F32").Value = 889
$ws.Range("F33").Value = 185
$ws.Range("F34").Value = 1871
$ws.Range("F39").Value = 148
$ws.Range("F40").Value = 73
$ws.Range("F41").Value = 28
$ws.Range("F43").Value = 86
$ws.Range("F46").Value = 79
$ws.Range("F47").Value = 37
$ws.Range("F48").Value = 1842
